# "Second commit on 3rd"
#
# 1. Login sheet: append a new credential row (admin2345 / admin5678)
#    formatted like the existing credential rows (rows 3-5, blue Consolas).
# 2. SubCategory sheet: the "link" column (B) data is no longer needed -
#    clear it out (rows 1-2), shrinking the sheet back to a single column.
# 3. Switch the active/selected tab from "Adminuser" to "SubCategory".

$wb = $excel.ActiveWorkbook

# --- Login sheet: add new admin credentials row ---------------------------
$wsLogin = $wb.Worksheets.Item("Login")

$wsLogin.Range("A6").Value = "admin2345"
$wsLogin.Range("B6").Value = "admin5678"

# Copy the formatting from the row above (A5:B5) onto the new row so it
# matches the existing admin-credential rows' style (blue Consolas font).
$wsLogin.Range("A5:B5").Copy()
$wsLogin.Range("A6:B6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$wsLogin.Range("H14").Select()

# --- SubCategory sheet: drop the link column data --------------------------
$wsSubCategory = $wb.Worksheets.Item("SubCategory")

$wsSubCategory.Range("B1:B2").ClearContents()

# --- Make SubCategory the active sheet/selection ---------------------------
$wsSubCategory.Activate()
$wsSubCategory.Range("B7").Select()
